$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 148.35
$ws.Range("I33").Value = 156.88235
$ws.Range("K33").Value = 156.88235
$ws.Range("M33").Value = 72.11765
$ws.Range("H111").Value = 1476.1666
$ws.Range("I111").Value = 609.6667
$ws.Range("J111").Value = 2342.6667
$ws.Range("K111").Value = 1829.0001
$ws.Range("L111").Value = 7028.000100000001
$ws.Range("M111").Value = 1237.9999
$ws.Range("N111").Value = -13162.0001
$ws.Range("H141").Value = 1832.7838
$ws.Range("I141").Value = 628.375
$ws.Range("K141").Value = 1885.125
$ws.Range("M141").Value = 3294.875

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 54764.895
$ws.Range("I2").Value = 85695.164
$ws.Range("K2").Value = 85695.164
$ws.Range("M2").Value = -85582.164
$ws.Range("H7").Value = 38500
$ws.Range("J7").Value = 38500
$ws.Range("L7").Value = 38500
$ws.Range("N7").Value = -38728
$ws.Range("H34").Value = 0
$ws.Range("J34").Value = 0
$ws.Range("L34").Value = 0
$ws.Range("N34").ClearContents()
$ws.Range("H45").Value = 1161.7
$ws.Range("I45").Value = 735.2222
$ws.Range("K45").Value = 735.2222
$ws.Range("M45").Value = -358.2222
$ws.Range("H110").Value = 624.875
$ws.Range("I110").Value = 624.875
$ws.Range("K110").Value = 624.875
$ws.Range("M110").Value = 1420.125
$ws.Range("H116").Value = 54764.895
$ws.Range("I116").Value = 85695.164
$ws.Range("K116").Value = 85695.164
$ws.Range("M116").Value = -83401.164
$ws.Range("H122").Value = 928.6
$ws.Range("I122").Value = 912.5714
$ws.Range("K122").Value = 2737.7142
$ws.Range("M122").Value = -287.7142000000003

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 54764.895
$ws.Range("I3").Value = 85695.164
$ws.Range("K3").Value = 85695.164
$ws.Range("M3").Value = -85581.164
$ws.Range("H107").Value = 29562.2
$ws.Range("I107").Value = 3337
$ws.Range("K107").Value = 3337
$ws.Range("M107").Value = -1417
$ws.Range("H134").Value = 73964.25
$ws.Range("I134").Value = 2616.8262
$ws.Range("J134").Value = 402162.4
$ws.Range("K134").Value = 7850.4786
$ws.Range("L134").Value = 1206487.2
$ws.Range("M134").Value = -5315.4786
$ws.Range("N134").Value = -1211557.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2311.449
$ws.Range("I31").Value = 2470.8057
$ws.Range("J31").Value = 1870.1538
$ws.Range("K31").Value = 2470.8057
$ws.Range("L31").Value = 1870.1538
$ws.Range("M31").Value = -2175.8057
$ws.Range("N31").Value = -2460.1538
$ws.Range("H34").Value = 2311.449
$ws.Range("I34").Value = 2470.8057
$ws.Range("J34").Value = 1870.1538
$ws.Range("K34").Value = 2470.8057
$ws.Range("L34").Value = 1870.1538
$ws.Range("M34").Value = -2268.8057
$ws.Range("N34").Value = -2274.1538
$ws.Range("H107").Value = 410.34784
$ws.Range("I107").Value = 442
$ws.Range("J107").Value = 403.6842
$ws.Range("K107").Value = 442
$ws.Range("L107").Value = 403.6842
$ws.Range("M107").Value = 1478
$ws.Range("N107").Value = -4243.6842
$ws.Range("H132").Value = 2098.24
$ws.Range("I132").Value = 1274
$ws.Range("J132").Value = 3849.75
$ws.Range("K132").Value = 3822
$ws.Range("L132").Value = 11549.25
$ws.Range("M132").Value = -1292
$ws.Range("N132").Value = -16609.25
$ws.Range("H141").Value = 50181.727
$ws.Range("J141").Value = 50181.727
$ws.Range("L141").Value = 50181.727
$ws.Range("N141").Value = -60541.727

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1633.7354
$ws.Range("I5").Value = 1410.2916
$ws.Range("J5").Value = 2170
$ws.Range("K5").Value = 4230.8748
$ws.Range("L5").Value = 6510
$ws.Range("M5").Value = -4118.8748
$ws.Range("N5").Value = -6734
$ws.Range("H129").Value = 73325.36
$ws.Range("J129").Value = 113296.11
$ws.Range("L129").Value = 339888.33
$ws.Range("N129").Value = -349888.33
$ws.Range("H135").Value = 1633.7354
$ws.Range("I135").Value = 1410.2916
$ws.Range("J135").Value = 2170
$ws.Range("K135").Value = 12692.6244
$ws.Range("L135").Value = 19530
$ws.Range("M135").Value = -10157.6244
$ws.Range("N135").Value = -24600

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H19").Value = 2000000
$ws.Range("I19").Value = 2000000
$ws.Range("J19").Value = 0
$ws.Range("K19").Value = 2000000
$ws.Range("L19").Value = 0
$ws.Range("M19").Value = -1999712
$ws.Range("N19").ClearContents()
$ws.Range("H74").Value = 69800
$ws.Range("J74").Value = 69800
$ws.Range("L74").Value = 69800
$ws.Range("N74").Value = -71672
$ws.Range("H77").Value = 69800
$ws.Range("J77").Value = 69800
$ws.Range("L77").Value = 209400
$ws.Range("N77").Value = -218760
$ws.Range("H107").Value = 621.6667
$ws.Range("I107").Value = 1550.5
$ws.Range("J107").Value = 157.25
$ws.Range("K107").Value = 1550.5
$ws.Range("L107").Value = 157.25
$ws.Range("M107").Value = 369.5
$ws.Range("N107").Value = -3997.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H29").Value = 12800
$ws.Range("J29").Value = 12800
$ws.Range("L29").Value = 12800
$ws.Range("N29").Value = -13390
$ws.Range("H61").Value = 1910.4445
$ws.Range("I61").Value = 1961.75
$ws.Range("J61").Value = 1500
$ws.Range("K61").Value = 1961.75
$ws.Range("L61").Value = 1500
$ws.Range("M61").Value = -1759.75
$ws.Range("N61").Value = -1904
$ws.Range("H113").Value = 1910.4445
$ws.Range("I113").Value = 1961.75
$ws.Range("J113").Value = 1500
$ws.Range("K113").Value = 1961.75
$ws.Range("L113").Value = 1500
$ws.Range("M113").Value = 208.25
$ws.Range("N113").Value = -5840
$ws.Range("H122").Value = 9671.643
$ws.Range("I122").Value = 10107.615
$ws.Range("J122").Value = 4004
$ws.Range("K122").Value = 30322.845
$ws.Range("L122").Value = 12012
$ws.Range("M122").Value = -27872.845
$ws.Range("N122").Value = -16912
$ws.Range("H136").Value = 2564.639
$ws.Range("I136").Value = 1661.08
$ws.Range("K136").Value = 4983.24
$ws.Range("M136").Value = -2433.24

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H16").Value = 49000
$ws.Range("J16").Value = 49000
$ws.Range("L16").Value = 49000
$ws.Range("N16").Value = -49584
$ws.Range("H119").Value = 2530000
$ws.Range("J119").Value = 2530000
$ws.Range("L119").Value = 2530000
$ws.Range("N119").Value = -2539676
$ws.Range("H136").Value = 1664.1731
$ws.Range("I136").Value = 1734.6578
$ws.Range("J136").Value = 1472.8572
$ws.Range("K136").Value = 5203.9734
$ws.Range("L136").Value = 4418.571599999999
$ws.Range("M136").Value = -2653.9734
$ws.Range("N136").Value = -9518.571599999999
